$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.593051075935364
$ws.Range("B1").Value = 5.044400691986084
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.771298170089722
$ws.Range("E1").Value = 1.69006884098053
